# "add bulk copy tool"
#
# Renames the original sheet and adds a second "Bulk Copy ICBC" sheet that
# holds an Input_Path / Onput_Path header block on top of a copy of the
# existing Producer Code / Producer Name table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ICBC E-Stamp Tool"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Bulk Copy ICBC"

# Bring over the Producer Code/Producer Name table (bold header row + the
# AA..EE / Alfa..Echo data rows) verbatim, formatting included, shifted
# down two rows to make room for the new path header block.
$ws1.Range("A2:B7").Copy($ws2.Range("A3"))

# Re-use the big header look (18pt bold label / 18pt value) from the
# original "Path" row for the two new Input_Path / Onput_Path rows.
$ws1.Range("A1:B1").Copy($ws2.Range("A1"))
$ws1.Range("A1:B1").Copy($ws2.Range("A2"))

$ws2.Range("A1").Value = "Input_Path"
$ws2.Range("A2").Value = "Onput_Path"
$ws2.Range("B1").Value = "C:\Users\<USERNAME>\Desktop\Old ICBC Copies"
$ws2.Range("B2").Value = "C:\Users\<USERNAME>\Desktop\New ICBC Copies"

$ws2.Rows.Item(1).RowHeight = 23.25
$ws2.Rows.Item(2).RowHeight = 23.25

$ws2.Columns.Item(1).ColumnWidth = 18

# Restore / set the view state of each sheet (multi-cell selection on the
# first sheet, single-cell selection on the new one).
$ws1.Range("A3:B7").Select()

$ws2.Range("E8").Select()
